# Resize the single-cell signature table: shrink its width (4059 -> 3300 dxa)
# and nudge its indent from the page margin (5911 -> 5926 dxa).
# Word COM exposes these dimensions in points, so convert dxa (twips) -> points
# by dividing by 20 (1 pt = 20 dxa).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# <w:tblW w:w="3300" .../> on <w:tblPr>
$t.PreferredWidth = 3300 / 20

# <w:tblInd w:w="5926" .../> on <w:tblPr>
$t.Rows.LeftIndent = 5926 / 20

# <w:gridCol w:w="3300"/> on <w:tblGrid> and <w:tcW w:w="3300" .../> on the cell
$t.Cell(1, 1).Width = 3300 / 20
